$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.152.73'
$ws.Range("E2").Value = '  +10.28%  '
$ws.Range("D3").Value = '1.678.63'
$ws.Range("E3").Value = '  +6.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3733'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.97%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3439'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.30'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +12.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.193'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07298'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.38'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.104'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.767'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.28%  '
$ws.Range("D16").Value = '1.679.45'
$ws.Range("E16").Value = '  +6.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001106'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9996'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06715'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '81.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.111'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.64%  '
$ws.Range("D24").Value = '24.145.01'
$ws.Range("E24").Value = '  +10.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.421'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.365'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -9.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.670'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.15'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.64%  '
$ws.Range("D30").Value = '1.859.74'
$ws.Range("E30").Value = '  +6.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '126.97'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.398'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.049'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9772'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.757'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08458'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '12.32'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.968'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06432'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.357'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02344'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.262'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.97%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.2122'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6181'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9989'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.805'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '13.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5963'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.040'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '126.46'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07169'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.38%  '
